$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("data")

# --- Add the new "metadata" worksheet, positioned right after "data" ---
$metaWs = $wb.Worksheets.Add($null, $dataWs)
$metaWs.Name = "metadata"

# Match sheet-level properties (outline summary flags, page margins) used by "data"
$metaWs.Outline.SummaryRow = 1
$metaWs.Outline.SummaryColumn = 1
$metaWs.PageSetup.LeftMargin = $dataWs.PageSetup.LeftMargin
$metaWs.PageSetup.RightMargin = $dataWs.PageSetup.RightMargin
$metaWs.PageSetup.TopMargin = $dataWs.PageSetup.TopMargin
$metaWs.PageSetup.BottomMargin = $dataWs.PageSetup.BottomMargin
$metaWs.PageSetup.HeaderMargin = $dataWs.PageSetup.HeaderMargin
$metaWs.PageSetup.FooterMargin = $dataWs.PageSetup.FooterMargin

# --- Header row (B1:G1) ---
$metaWs.Range("B1").Value = "data_name"
$metaWs.Range("C1").Value = "data_id"
$metaWs.Range("D1").Value = "data_version"
$metaWs.Range("E1").Value = "data_version_created"
$metaWs.Range("F1").Value = "panel_query_time"
$metaWs.Range("G1").Value = "panel_get_request"

# Match header styling (bold, bordered, centered) used by the "data" sheet header
$dataWs.Range("B1").Copy()
$metaWs.Range("B1:G1").PasteSpecial(-4122)

# --- Data row (row 2) ---
$metaWs.Range("A2").Value = 0
$dataWs.Range("A2").Copy()
$metaWs.Range("A2").PasteSpecial(-4122)

$metaWs.Range("B2").Value = "Pituitary hormone deficiency"
$metaWs.Range("C2").Value = 483
$metaWs.Range("D2").NumberFormat = "@"
$metaWs.Range("D2").Value = "2.7"
$metaWs.Range("E2").Value = "2021-09-15T11:19:48.513650Z"
$metaWs.Range("F2").Value = "2021-10-05 14:22:13.469399"
$metaWs.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/483/?format=json"

# --- Refresh "data" sheet query timestamps (column F, rows 2-53) ---
$dataWs.Range("F2").Value = "2021-10-05 14:22:13.472869"
$dataWs.Range("F3").Value = "2021-10-05 14:22:13.472880"
$dataWs.Range("F4").Value = "2021-10-05 14:22:13.472884"
$dataWs.Range("F5").Value = "2021-10-05 14:22:13.472886"
$dataWs.Range("F6").Value = "2021-10-05 14:22:13.472889"
$dataWs.Range("F7").Value = "2021-10-05 14:22:13.472892"
$dataWs.Range("F8").Value = "2021-10-05 14:22:13.472895"
$dataWs.Range("F9").Value = "2021-10-05 14:22:13.472898"
$dataWs.Range("F10").Value = "2021-10-05 14:22:13.472901"
$dataWs.Range("F11").Value = "2021-10-05 14:22:13.472904"
$dataWs.Range("F12").Value = "2021-10-05 14:22:13.472907"
$dataWs.Range("F13").Value = "2021-10-05 14:22:13.472910"
$dataWs.Range("F14").Value = "2021-10-05 14:22:13.472912"
$dataWs.Range("F15").Value = "2021-10-05 14:22:13.472915"
$dataWs.Range("F16").Value = "2021-10-05 14:22:13.472917"
$dataWs.Range("F17").Value = "2021-10-05 14:22:13.472920"
$dataWs.Range("F18").Value = "2021-10-05 14:22:13.472923"
$dataWs.Range("F19").Value = "2021-10-05 14:22:13.472926"
$dataWs.Range("F20").Value = "2021-10-05 14:22:13.472929"
$dataWs.Range("F21").Value = "2021-10-05 14:22:13.472931"
$dataWs.Range("F22").Value = "2021-10-05 14:22:13.472934"
$dataWs.Range("F23").Value = "2021-10-05 14:22:13.472936"
$dataWs.Range("F24").Value = "2021-10-05 14:22:13.472939"
$dataWs.Range("F25").Value = "2021-10-05 14:22:13.472942"
$dataWs.Range("F26").Value = "2021-10-05 14:22:13.472945"
$dataWs.Range("F27").Value = "2021-10-05 14:22:13.472948"
$dataWs.Range("F28").Value = "2021-10-05 14:22:13.472950"
$dataWs.Range("F29").Value = "2021-10-05 14:22:13.472953"
$dataWs.Range("F30").Value = "2021-10-05 14:22:13.472956"
$dataWs.Range("F31").Value = "2021-10-05 14:22:13.472959"
$dataWs.Range("F32").Value = "2021-10-05 14:22:13.472962"
$dataWs.Range("F33").Value = "2021-10-05 14:22:13.472965"
$dataWs.Range("F34").Value = "2021-10-05 14:22:13.472968"
$dataWs.Range("F35").Value = "2021-10-05 14:22:13.472971"
$dataWs.Range("F36").Value = "2021-10-05 14:22:13.472974"
$dataWs.Range("F37").Value = "2021-10-05 14:22:13.472976"
$dataWs.Range("F38").Value = "2021-10-05 14:22:13.472979"
$dataWs.Range("F39").Value = "2021-10-05 14:22:13.472982"
$dataWs.Range("F40").Value = "2021-10-05 14:22:13.472984"
$dataWs.Range("F41").Value = "2021-10-05 14:22:13.472987"
$dataWs.Range("F42").Value = "2021-10-05 14:22:13.472990"
$dataWs.Range("F43").Value = "2021-10-05 14:22:13.472993"
$dataWs.Range("F44").Value = "2021-10-05 14:22:13.472996"
$dataWs.Range("F45").Value = "2021-10-05 14:22:13.472999"
$dataWs.Range("F46").Value = "2021-10-05 14:22:13.473001"
$dataWs.Range("F47").Value = "2021-10-05 14:22:13.473004"
$dataWs.Range("F48").Value = "2021-10-05 14:22:13.473007"
$dataWs.Range("F49").Value = "2021-10-05 14:22:13.473009"
$dataWs.Range("F50").Value = "2021-10-05 14:22:13.473012"
$dataWs.Range("F51").Value = "2021-10-05 14:22:13.473015"
$dataWs.Range("F52").Value = "2021-10-05 14:22:13.473018"
$dataWs.Range("F53").Value = "2021-10-05 14:22:13.473020"
